$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 92
$ws.Range("D92").Value = 44782
$ws.Range("D92").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K92").Value = 40000
$ws.Range("L92").Value = 42000
$ws.Range("M92").Value = 41000
$ws.Range("P92").Value = 586

# Row 93
$ws.Range("D93").Value = 44427
$ws.Range("D93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H93").Value = 'Dulce o Americano'
$ws.Range("J93").Value = 100
$ws.Range("K93").Value = 34000
$ws.Range("L93").Value = 35000
$ws.Range("M93").Value = 34500
$ws.Range("N93").Value = '$/malla 70 unidades'
$ws.Range("O93").Value = 'Región de Arica y Parinacota'
$ws.Range("P93").Value = 493
$ws.Range("Q93").Value = 70

# Row 94
$ws.Range("I94").Value = 'Primera'
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 200
$ws.Range("L94").Value = 200
$ws.Range("M94").Value = 200
$ws.Range("P94").Value = 200

# Row 95
$ws.Range("D95").Value = 44628
$ws.Range("D95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H95").Value = 'Choclero'
$ws.Range("I95").Value = 'Segunda'
$ws.Range("J95").Value = 3500
$ws.Range("K95").Value = 150
$ws.Range("L95").Value = 150
$ws.Range("M95").Value = 150
$ws.Range("N95").Value = '$/unidad'
$ws.Range("O95").Value = 'Región de O''Higgins'
$ws.Range("P95").Value = 150
$ws.Range("Q95").Value = 1

# Row 96
$ws.Range("D96").Value = 44483
$ws.Range("D96").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H96").Value = 'Dulce o Americano'
$ws.Range("J96").Value = 450
$ws.Range("K96").Value = 25000
$ws.Range("L96").Value = 26000
$ws.Range("M96").Value = 25556
$ws.Range("N96").Value = '$/malla 70 unidades'
$ws.Range("O96").Value = 'Región de Arica y Parinacota'
$ws.Range("P96").Value = 365
$ws.Range("Q96").Value = 70

# Row 97
$ws.Range("I97").Value = 'Primera'
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 300
$ws.Range("L97").Value = 350
$ws.Range("M97").Value = 325
$ws.Range("P97").Value = 325

# Row 98
$ws.Range("D98").Value = 44204
$ws.Range("D98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I98").Value = 'Segunda'
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 250
$ws.Range("L98").Value = 250
$ws.Range("M98").Value = 250
$ws.Range("O98").Value = 'Región Metropolitana'
$ws.Range("P98").Value = 250

# Row 99
$ws.Range("D99").Value = 44484
$ws.Range("D99").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H99").Value = 'Dulce o Americano'
$ws.Range("I99").Value = 'Primera'
$ws.Range("J99").Value = 150
$ws.Range("K99").Value = 43000
$ws.Range("L99").Value = 45000
$ws.Range("M99").Value = 44067
$ws.Range("N99").Value = '$/malla 70 unidades'
$ws.Range("O99").Value = 'Región de Arica y Parinacota'
$ws.Range("P99").Value = 630
$ws.Range("Q99").Value = 70

# Row 100
$ws.Range("D100").Value = 44238
$ws.Range("D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 300
$ws.Range("L100").Value = 320
$ws.Range("M100").Value = 310
$ws.Range("P100").Value = 310

# Row 101
$ws.Range("D101").Value = 44238
$ws.Range("D101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K101").Value = 250
$ws.Range("L101").Value = 250
$ws.Range("M101").Value = 250
$ws.Range("P101").Value = 250

# Row 102
$ws.Range("D102").Value = 44279
$ws.Range("D102").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H102").Value = 'Choclero'
$ws.Range("J102").Value = 20000
$ws.Range("K102").Value = 230
$ws.Range("L102").Value = 250
$ws.Range("M102").Value = 240
$ws.Range("N102").Value = '$/unidad'
$ws.Range("O102").Value = 'Región de O''Higgins'
$ws.Range("P102").Value = 240
$ws.Range("Q102").Value = 1

# Row 103
$ws.Range("D103").Value = 44279
$ws.Range("D103").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I103").Value = 'Segunda'
$ws.Range("J103").Value = 5000
$ws.Range("L103").Value = 200
$ws.Range("M103").Value = 200
$ws.Range("O103").Value = 'Región de O''Higgins'
$ws.Range("P103").Value = 200

# Row 104
$ws.Range("D104").Value = 44355
$ws.Range("D104").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J104").Value = 50
$ws.Range("K104").Value = 11000
$ws.Range("L104").Value = 12000
$ws.Range("M104").Value = 11400
$ws.Range("N104").Value = '$/malla 60 unidades'
$ws.Range("O104").Value = 'Provincia de Limarí'
$ws.Range("P104").Value = 190
$ws.Range("Q104").Value = 60

# Row 105
$ws.Range("A105").Value = 11
$ws.Range("B105").Value = 'Vega Monumental Concepción'
$ws.Range("C105").Value = 'Bíobío'
$ws.Range("D105").Value = 44657
$ws.Range("D105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E105").Value = 8
$ws.Range("F105").Value = 100112024
$ws.Range("G105").Value = 'Choclo'
$ws.Range("H105").Value = 'Choclero'
$ws.Range("I105").Value = 'Primera'
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 200
$ws.Range("L105").Value = 250
$ws.Range("M105").Value = 225
$ws.Range("N105").Value = '$/unidad'
$ws.Range("O105").Value = 'Región Metropolitana'
$ws.Range("P105").Value = 225
$ws.Range("Q105").Value = 1
$ws.Range("R105").Value = 'Hortaliza'

# Row 106
$ws.Range("A106").Value = 11
$ws.Range("B106").Value = 'Vega Monumental Concepción'
$ws.Range("C106").Value = 'Bíobío'
$ws.Range("D106").Value = 44453
$ws.Range("D106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E106").Value = 8
$ws.Range("F106").Value = 100112024
$ws.Range("G106").Value = 'Choclo'
$ws.Range("H106").Value = 'Dulce o Americano'
$ws.Range("I106").Value = 'Primera'
$ws.Range("J106").Value = 100
$ws.Range("K106").Value = 36000
$ws.Range("L106").Value = 38000
$ws.Range("M106").Value = 37000
$ws.Range("N106").Value = '$/malla 70 unidades'
$ws.Range("O106").Value = 'Región de Arica y Parinacota'
$ws.Range("P106").Value = 529
$ws.Range("Q106").Value = 70
$ws.Range("R106").Value = 'Hortaliza'
